$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("B2:B6")
$range.NumberFormat = "@"

$ws.Range("B2").Value = "07947112241"
$ws.Range("B3").Value = "07947435158"
$ws.Range("B4").Value = "07947124072"
$ws.Range("B5").Value = "07942698613"
$ws.Range("B6").Value = "07947137139"

$range.ClearFormats()
